# Fill in the previously-blank "Notes to Accounts" amount columns (C/D) on
# the Balance Sheet, and correct the FY24 current-assets total in C34.
#
# These cells store plain numbers as literal text (e.g. "2,046.90") rather
# than numeric values with a number format, matching the rest of the sheet.
# A leading apostrophe forces Excel to keep the entry as text instead of
# auto-converting it to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value  = "'204.86"
$ws.Range("D6").Value  = "'2,046.90"

$ws.Range("D7").Value  = "'9,958.15"

$ws.Range("C8").Value  = "'204.86"
$ws.Range("D8").Value  = "'12,005.06"

$ws.Range("D10").Value = "'914.46"

$ws.Range("D11").Value = "'49.00"

$ws.Range("D12").Value = "'963.46"

$ws.Range("C14").Value = "'491.39"
$ws.Range("D14").Value = "'104.23"

$ws.Range("D15").Value = "'125.35"

$ws.Range("D16").Value = "'179.27"

$ws.Range("C17").Value = "'491.39"
$ws.Range("D17").Value = "'408.85"

$ws.Range("C18").Value = "'696.25"
$ws.Range("D18").Value = "'13,377.37"

$ws.Range("C25").Value = "'0.60"

$ws.Range("C26").Value = "'0.60"

# FY24 total current assets corrected from 3,884.46 to 3,885.06
$ws.Range("C34").Value = "'3,885.06"
